$wb = $excel.ActiveWorkbook

$wsStories  = $wb.Worksheets.Item("Stories")
$wsBurndown = $wb.Worksheets.Item("Burn down")

# --- Stories sheet: F3:F5 flip from "Ja" to "Nein" ---
$wsStories.Range("F3").Value = "Nein"
$wsStories.Range("F4").Value = "Nein"
$wsStories.Range("F5").Value = "Nein"

# --- Stories sheet: fill in the "Zustaendig" (team) column for rows 6-10 ---
$wsStories.Range("C6").Value  = "Team2"
$wsStories.Range("C7").Value  = "Team3"
$wsStories.Range("C8").Value  = "Team1"
$wsStories.Range("C9").Value  = "Team2"
$wsStories.Range("C10").Value = "Team3"

# --- Burn down sheet: record the actual (Ist) burn down value for iteration 3 ---
$wsBurndown.Range("B4").Value = 34

# --- View / selection state: Burn down was active before, Stories becomes active ---
# Select on the sheet that should end up NOT active first ...
[void]$wsBurndown.Range("B5").Select()

# ... then activate + select on the sheet that should end up active last,
# so it "wins" the workbook-level activeTab / tabSelected state.
[void]$wsStories.Activate()
[void]$wsStories.Range("F5").Select()
